$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.535.40"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.438.45"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "594.21"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "135.44"
$ws.Range("E6").Value = "  -8.44%  "
$ws.Range("D7").Value = "3.439.05"
$ws.Range("E7").Value = "  -4.23%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "7.55"
$ws.Range("E10").Value = "  -3.57%  "
$ws.Range("E11").Value = "  -10.12%  "
$ws.Range("D12").Value = "0.379"
$ws.Range("E12").Value = "  -8.25%  "
$ws.Range("D13").Value = "4.017.24"
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("E14").Value = "  -12.39%  "
$ws.Range("D15").Value = "26.44"
$ws.Range("E15").Value = "  -10.44%  "
$ws.Range("D16").Value = "3.447.60"
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "65.432.11"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "0.114"
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "10.04"
$ws.Range("E19").Value = "  -9.01%  "
$ws.Range("D20").Value = "5.76"
$ws.Range("E20").Value = "  -8.95%  "
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -7.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.32%  "
$ws.Range("E23").Value = "  -10.55%  "
$ws.Range("D24").Value = "73.21"
$ws.Range("E24").Value = "  -6.67%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "3.579.15"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -12.61%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -11.35%  "
$ws.Range("E30").Value = "  -9.14%  "
$ws.Range("D31").Value = "8.15"
$ws.Range("E31").Value = "  -12.83%  "
$ws.Range("D32").Value = "3.444.42"
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("D35").Value = "22.69"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.20"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("E37").Value = "  -13.83%  "
$ws.Range("D38").Value = "6.86"
$ws.Range("E38").Value = "  -11.42%  "
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("E40").Value = "  -13.62%  "
$ws.Range("D41").Value = "0.0773"
$ws.Range("E41").Value = "  -9.31%  "
$ws.Range("E42").Value = "  -7.84%  "
$ws.Range("D43").Value = "43.51"
$ws.Range("E43").Value = "  -5.21%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "4.41"
$ws.Range("E45").Value = "  -14.98%  "
$ws.Range("E46").Value = "  -12.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.70"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("E49").Value = "  -8.70%  "
$ws.Range("E50").Value = "  -16.27%  "
$ws.Range("D51").Value = "2.189.43"
$ws.Range("E51").Value = "  -8.29%  "
